$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1030.3939
$ws.Range("I17").Value = 580.5862
$ws.Range("K17").Value = 1741.7586
$ws.Range("M17").Value = -1573.7586

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2923.4707
$ws.Range("I64").Value = 2842.8572
$ws.Range("J64").Value = 2979.9
$ws.Range("K64").Value = 2842.8572
$ws.Range("L64").Value = 2979.9
$ws.Range("M64").Value = -2594.8572
$ws.Range("N64").Value = -3475.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2923.4707
$ws.Range("I67").Value = 2842.8572
$ws.Range("J67").Value = 2979.9
$ws.Range("K67").Value = 2842.8572
$ws.Range("L67").Value = 2979.9
$ws.Range("M67").Value = -1984.8572
$ws.Range("N67").Value = -4695.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3204
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 3204
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 3204
$ws.Range("N74").Value = -5076
$ws.Range("M74").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4764865
$ws.Range("I76").Value = 5211290
$ws.Range("K76").Value = 5211290
$ws.Range("M76").Value = -5210975

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3204
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 3204
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 16020
$ws.Range("N77").Value = -25380
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4764865
$ws.Range("I79").Value = 5211290
$ws.Range("K79").Value = 5211290
$ws.Range("M79").Value = -5210198

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 33251.066
$ws.Range("I63").Value = 78504.836
$ws.Range("J63").Value = 4669.737
$ws.Range("K63").Value = 78504.836
$ws.Range("L63").Value = 4669.737
$ws.Range("M63").Value = -77818.836
$ws.Range("N63").Value = -6041.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 33251.066
$ws.Range("I66").Value = 78504.836
$ws.Range("J66").Value = 4669.737
$ws.Range("K66").Value = 392524.18
$ws.Range("L66").Value = 23348.685
$ws.Range("M66").Value = -389092.18
$ws.Range("N66").Value = -30212.685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 50377.715
$ws.Range("I122").Value = 92593.82000000001
$ws.Range("J122").Value = 3940
$ws.Range("K122").Value = 277781.46
$ws.Range("L122").Value = 11820
$ws.Range("M122").Value = -275331.46
$ws.Range("N122").Value = -16720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1591.2444
$ws.Range("I86").Value = 1456.7949
$ws.Range("J86").Value = 2465.1667
$ws.Range("K86").Value = 1456.7949
$ws.Range("L86").Value = 2465.1667
$ws.Range("M86").Value = -333.7949000000001
$ws.Range("N86").Value = -4711.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1591.2444
$ws.Range("I89").Value = 1456.7949
$ws.Range("J89").Value = 2465.1667
$ws.Range("K89").Value = 7283.9745
$ws.Range("L89").Value = 12325.8335
$ws.Range("M89").Value = -1667.9745
$ws.Range("N89").Value = -23557.8335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 167666.67
$ws.Range("I94").Value = 250775
$ws.Range("J94").Value = 1450
$ws.Range("K94").Value = 250775
$ws.Range("L94").Value = 1450
$ws.Range("M94").Value = -250324
$ws.Range("N94").Value = -2352

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4210
$ws.Range("I62").Value = 4237.5
$ws.Range("J62").Value = 4100
$ws.Range("K62").Value = 4237.5
$ws.Range("L62").Value = 4100
$ws.Range("M62").Value = -3613.5
$ws.Range("N62").Value = -5348

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4210
$ws.Range("I65").Value = 4237.5
$ws.Range("J65").Value = 4100
$ws.Range("K65").Value = 21187.5
$ws.Range("L65").Value = 20500
$ws.Range("M65").Value = -18067.5
$ws.Range("N65").Value = -26740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 668.1836499999999
$ws.Range("I5").Value = 489.2647
$ws.Range("J5").Value = 1073.7333
$ws.Range("K5").Value = 1467.7941
$ws.Range("L5").Value = 3221.199900000001
$ws.Range("M5").Value = -1355.7941
$ws.Range("N5").Value = -3445.199900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 864
$ws.Range("I22").Value = 575
$ws.Range("K22").Value = 1725
$ws.Range("M22").Value = -1556

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 463.10526
$ws.Range("I25").Value = 100
$ws.Range("J25").Value = 483.27777
$ws.Range("K25").Value = 300
$ws.Range("L25").Value = 1449.83331
$ws.Range("M25").Value = -131
$ws.Range("N25").Value = -1787.83331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 864
$ws.Range("I27").Value = 575
$ws.Range("K27").Value = 1725
$ws.Range("M27").Value = -1623

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 19608142
$ws.Range("I29").Value = 69
$ws.Range("J29").Value = 33333794
$ws.Range("K29").Value = 207
$ws.Range("L29").Value = 100001382
$ws.Range("M29").Value = 70
$ws.Range("N29").Value = -100001936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 463.10526
$ws.Range("I30").Value = 100
$ws.Range("J30").Value = 483.27777
$ws.Range("K30").Value = 300
$ws.Range("L30").Value = 1449.83331
$ws.Range("M30").Value = -198
$ws.Range("N30").Value = -1653.83331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 626.2
$ws.Range("I113").Value = 579.36365
$ws.Range("J113").Value = 755
$ws.Range("K113").Value = 1738.09095
$ws.Range("L113").Value = 2265
$ws.Range("M113").Value = 431.90905
$ws.Range("N113").Value = -6605

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 668.1836499999999
$ws.Range("I135").Value = 489.2647
$ws.Range("J135").Value = 1073.7333
$ws.Range("K135").Value = 4403.3823
$ws.Range("L135").Value = 9663.599700000001
$ws.Range("M135").Value = -1868.3823
$ws.Range("N135").Value = -14733.5997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8492.134
$ws.Range("I70").Value = 8898.615
$ws.Range("K70").Value = 8898.615
$ws.Range("M70").Value = -8628.615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8492.134
$ws.Range("I73").Value = 8898.615
$ws.Range("K73").Value = 8898.615
$ws.Range("M73").Value = -7962.615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 856849.25
$ws.Range("I80").Value = 2252045
$ws.Range("J80").Value = 59594.57
$ws.Range("K80").Value = 2252045
$ws.Range("L80").Value = 59594.57
$ws.Range("M80").Value = -2251047
$ws.Range("N80").Value = -61590.57

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 856849.25
$ws.Range("I83").Value = 2252045
$ws.Range("J83").Value = 59594.57
$ws.Range("K83").Value = 11260225
$ws.Range("L83").Value = 297972.85
$ws.Range("M83").Value = -11255233
$ws.Range("N83").Value = -307956.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 90009.336
$ws.Range("J43").Value = 90009.336
$ws.Range("L43").Value = 90009.336
$ws.Range("N43").Value = -90395.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6208165.5
$ws.Range("I132").Value = 2591.2856
$ws.Range("K132").Value = 7773.8568
$ws.Range("M132").Value = -5243.8568
